$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56

$ws.Cells.Item($row, 1).Value = "2025-10-02 21:24:40"
$ws.Cells.Item($row, 2).Value = "Noah"
$ws.Cells.Item($row, 3).Value = 8450689526

# D56 must be stored as text ("13052054965"), not a number, matching most
# other rows in the Phone column. A leading apostrophe forces text entry so
# the numeric-looking string isn't coerced to a number; reapplying the
# "Normal" style afterwards drops the resulting quote-prefix formatting so
# the cell ends up plain/unstyled like its siblings.
$ws.Cells.Item($row, 4).Value = "'13052054965"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = " "
$ws.Cells.Item($row, 6).Value = "my-node-server/public/uploads/images\photo_2025-10-03_01-24-41.jpg"

# G56 must exist as an empty (but present) cell, matching the sheet's
# existing "empty inline string" cells (e.g. F2, G2, F55, G55). Writing ""
# directly clears/omits the cell entirely, so force text entry via a
# leading apostrophe, then drop the resulting quote-prefix style so the
# cell matches the plain, unstyled empty cells elsewhere in the sheet.
$ws.Cells.Item($row, 7).Value = "'"
$ws.Cells.Item($row, 7).Style = "Normal"
